# Add files via upload
# Sujeto_1/Carbohidrates.xlsx — fix a handful of "Value (g)" cells that
# had been entered/saved as raw numbers (with a "#,##0" style) instead of
# the plain-text decimal strings used everywhere else in column B.
#
#   B8  : 39375 -> "39.38"
#   B18 : 20305 -> "20.3"
#   B20:B27 : 30875 -> "30.88" (all eight rows share the same text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force text entry (like Excel does when a cell is pre-formatted as
    # Text) so the numeric-looking string isn't re-interpreted as a
    # number, then drop back to the workbook's default "Normal" style so
    # no leftover number format sticks around on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("B8") "39.38"
Set-TextValue $ws.Range("B18") "20.3"
Set-TextValue $ws.Range("B20:B27") "30.88"

# The saved file also scrolled the sheet view down to row 24 (was row 16)
# and dropped the lingering D19 selection / column B width override —
# cosmetic leftovers from the author's last interaction before saving.
$ws.Range("A24").Select() | Out-Null
$ws.Columns("B:B").AutoFit() | Out-Null

Write-Output "done"
